$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label row (row 16) - new env description string
$ws.Range("A16").Value = "Env: Multiple threads for reading (thread count=4) and single thread for writes. Version 0.0.3"

# Raw timing data for rows 17-20 (columns A,B for the left block; G,H for the right block)
$ws.Range("A17").Value = 0.4828587962962963
$ws.Range("B17").Value = 0.48312500000000003
$ws.Range("G17").Value = 0.48372685185185182
$ws.Range("H17").Value = 0.48376157407407411

$ws.Range("A18").Value = 0.48431712962962964
$ws.Range("B18").Value = 0.48462962962962958
$ws.Range("G18").Value = 0.48494212962962963
$ws.Range("H18").Value = 0.48497685185185185

$ws.Range("A19").Value = 0.48835648148148153
$ws.Range("B19").Value = 0.48863425925925924
$ws.Range("G19").Value = 0.48885416666666665
$ws.Range("H19").Value = 0.48887731481481483

$ws.Range("A20").Value = 0.48984953703703704
$ws.Range("B20").Value = 0.49017361111111107
$ws.Range("G20").Value = 0.49032407407407402
$ws.Range("H20").Value = 0.49034722222222221

# Apply the time number format (matches style index 1 used elsewhere in the sheet)
$ws.Range("A17:B20").NumberFormat = "h:mm:ss"
$ws.Range("G17:H20").NumberFormat = "h:mm:ss"

# Difference formulas per row
$ws.Range("C17").Formula = "=B17-A17"
$ws.Range("I17").Formula = "=H17-G17"
$ws.Range("C18").Formula = "=B18-A18"
$ws.Range("I18").Formula = "=H18-G18"
$ws.Range("C19").Formula = "=B19-A19"
$ws.Range("I19").Formula = "=H19-G19"
$ws.Range("C20").Formula = "=B20-A20"
$ws.Range("I20").Formula = "=H20-G20"

$ws.Range("C17:C20").NumberFormat = "h:mm:ss"
$ws.Range("I17:I20").NumberFormat = "h:mm:ss"

# Averages in row 21
$ws.Range("C21").Formula = "=AVERAGE(C17:C20)"
$ws.Range("I21").Formula = "=AVERAGE(I17:I20)"
$ws.Range("C21").NumberFormat = "h:mm:ss"
$ws.Range("I21").NumberFormat = "h:mm:ss"

# Update the selection to match the author's final cursor position
$ws.Range("I23").Select()
